# Se procesan de nuevo los datos con las nuevas dimensiones curadas
#
# This reprocesses the "provincia" column metadata block (column E, rows 2-4):
#   - E2: sdmx-dimension:refArea  -> iaest-measure:provincia   (it's now a measure, not a dimension)
#   - E3: dim                    -> medida                     (dim/medida marker follows suit)
#   - E4: URI-Provincia          -> xsd:int                    (datatype changes from URI to int)
#
# The now-unused shared string "URI-Provincia" is dropped automatically by Excel
# when no cell references it anymore.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "iaest-measure:provincia"
$ws.Range("E3").Value = "medida"
$ws.Range("E4").Value = "xsd:int"
